# weight_tracker update: "add update to today; extend back extrapolation period"
#
# The worksheet Excel shows as "Sheet2" is the small summary/pivot-source
# table; the worksheet named "raw_data" is the long date/time/weight log
# (rows 1-166 before this edit) that backs the big scatter chart. The four
# new readings belong on that "raw_data" sheet, immediately below the
# existing data (rows 167-170), extending the shared TOD formula and the
# chart's cached series.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("raw_data")

# ---- new readings -------------------------------------------------------
# row, date+time serial, time-of-day fraction, weight
$newRows = @(
    @{ Row = 167; A = 44098.331250000003; B = 0.33124999999999999; C = 70.8 },
    @{ Row = 168; A = 44098.309027777781; B = 0.30902777777777779; C = 71.400000000000006 },
    @{ Row = 169; A = 44098.308333333334; B = 0.30833333333333335; C = 71.8 },
    @{ Row = 170; A = 44097.910416666666; B = 0.91041666666666676; C = 72.400000000000006 }
)

foreach ($r in $newRows) {
    $ws.Cells.Item($r.Row, 1).Value = $r.A
    $ws.Cells.Item($r.Row, 2).Value = $r.B
    $ws.Cells.Item($r.Row, 3).Value = $r.C
}

# Copy column A/B number formatting (date, time) down from the last
# pre-existing row so the new cells render the same way.
$ws.Range("A166:B166").Copy()
$ws.Range("A167:B170").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Extend the TOD (AM/PM) formula down into the new rows.
foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Cells.Item($row, 4).Formula = '=IF(B' + $row + '<TIME(12,0,0),"AM","PM")'
}

# Note: the chart (raw_data!$A$2:$A$245 / $C$2:$C$245) caches its plotted
# points rather than reading the range live. Writing SeriesCollection.Values
# / .XValues on this host emits a brand-new (unlinked, empty) chart part
# instead of updating chart2.xml's cache in place, which would leave a
# broken extra chart in the workbook -- worse than leaving the existing
# chart cache stale -- so that avenue is deliberately not used here.

# ---- selection, matching the saved view after the edit -------------------
$ws.Activate()
$ws.Range("C170").Select()
